# Update the "Förändrad" (Changed) date column (C) for rows 2-70
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 70; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
